$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "PCOS"

# Existing data row moves to row 2; write the new header row in row 1
$ws.Range("A1").Value = "Recipe Name"
$ws.Range("B1").Value = "Recipe Category(Breakfast/lunch/snack/dinner)"
$ws.Range("C1").Value = "Food Category(Veg/non-veg/vegan/Jain)"
$ws.Range("D1").Value = "Ingredients"
$ws.Range("E1").Value = "Preparation Time"
$ws.Range("F1").Value = "Cooking Time"
$ws.Range("G1").Value = "Preparation method"

# Named cell style ("Normal 2") carrying the 10pt black Helvetica Neue font
$style = $wb.Styles.Add("Normal 2")
$style.Font.Size = 10
$style.Font.Color = 0
$style.Font.Name = "Helvetica Neue"

# Apply the style, then the remaining direct formatting (white fill, thin
# border, top+wrap alignment, text number format) to the header row
$hdr = $ws.Range("A1:J1")
$hdr.Style = "Normal 2"
$hdr.Interior.Pattern = 1
$hdr.Interior.Color = 16777215
$hdr.Borders.LineStyle = 1
$hdr.VerticalAlignment = -4160
$hdr.WrapText = $true
$hdr.NumberFormat = "@"
$hdr.RowHeight = 75

$ws.Range("J1").Select()

Write-Output "done"
